$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Database update: drop the oldest "12 ماهه منتهی به 1396/12" period column,
# shift the remaining four periods one column to the left (D<-E<-F<-G<-H) and
# append the newest "12 ماهه منتهی به 1401/12" period as the new column H.
# Same shift for the "تاریخ انتشار" (publish date) row, with a brand new
# publish-date pair for the new column.
# ---------------------------------------------------------------------------

# Row 8 - financial period headers
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# Row 9 - publish dates
$ws.Range("D9").Value = "1399-02-14 (10)"
$ws.Range("E9").Value = "1400-02-29 (9)"
$ws.Range("F9").Value = "1401-04-04 (9)"
$ws.Range("G9").Value = "1402-02-27 (9)"
$ws.Range("H9").Value = "1402-02-27 (2)"

# Row 11 - فروش (Sales)
$ws.Range("D11").Value = 101195
$ws.Range("E11").Value = 109441
$ws.Range("F11").Value = 106212
$ws.Range("G11").Value = 231624
$ws.Range("H11").Value = 224666

# Row 12 - بهای تمام شده کالای فروش رفته (Cost of goods sold)
$ws.Range("D12").Value = -37842
$ws.Range("E12").Value = -44577
$ws.Range("F12").Value = -33629
$ws.Range("G12").Value = -81315
$ws.Range("H12").Value = -72185

# Row 13 - سود (زیان) ناخالص (Gross profit)
$ws.Range("D13").Value = 63353
$ws.Range("E13").Value = 64864
$ws.Range("F13").Value = 72583
$ws.Range("G13").Value = 150308
$ws.Range("H13").Value = 152481

# Row 14 - هزینه های عمومی, اداری و تشکیلاتی (G&A expenses)
$ws.Range("D14").Value = -7922
$ws.Range("E14").Value = -8428
$ws.Range("F14").Value = -6841
$ws.Range("G14").Value = -8562
$ws.Range("H14").Value = -10670

# Row 15 - هزینه کاهش ارزش دریافتنی‌‏ها : unchanged ("-" in every column)

# Row 16 - خالص سایر درامدها (هزینه ها) ی عملیاتی
$ws.Range("D16").Value = -1268
$ws.Range("E16").Value = 1096
$ws.Range("F16").Value = 123
$ws.Range("G16").Value = 1439
$ws.Range("H16").Value = 9823

# Row 17 - سود (زیان) عملیاتی (Operating profit)
$ws.Range("D17").Value = 54163
$ws.Range("E17").Value = 57532
$ws.Range("F17").Value = 65865
$ws.Range("G17").Value = 143185
$ws.Range("H17").Value = 151634

# Row 18 - هزینه های مالی (Financial expenses)
$ws.Range("D18").Value = -3775
$ws.Range("E18").Value = -1676
$ws.Range("F18").Value = -769
$ws.Range("G18").Value = -971
$ws.Range("H18").Value = -472

# Row 19 - خالص سایر درامدها و هزینه های غیرعملیاتی
$ws.Range("D19").Value = 399
$ws.Range("E19").Value = -2140
$ws.Range("F19").Value = -10109
$ws.Range("G19").Value = -3734
$ws.Range("H19").Value = -19161

# Row 20 - سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
$ws.Range("D20").Value = 50787
$ws.Range("E20").Value = 53716
$ws.Range("F20").Value = 54987
$ws.Range("G20").Value = 138479
$ws.Range("H20").Value = 132001

# Row 21 - مالیات (Tax): moves from column H to column G, column H becomes "-"
$ws.Range("G21").Value = 10396
$ws.Range("H21").Value = "-"

# Row 22 - سود (زیان) خالص عملیات در حال تداوم
$ws.Range("D22").Value = 50787
$ws.Range("E22").Value = 53716
$ws.Range("F22").Value = 54987
$ws.Range("G22").Value = 148875
$ws.Range("H22").Value = 132001

# Row 23 - سود (زیان) عملیات متوقف شده پس از اثر مالیاتی : unchanged ("-")

# Row 24 - سود (زیان) خالص (Net profit)
$ws.Range("D24").Value = 50787
$ws.Range("E24").Value = 53716
$ws.Range("F24").Value = 54987
$ws.Range("G24").Value = 148875
$ws.Range("H24").Value = 132001

# Row 25 - سود هر سهم پس از کسر مالیات : unchanged (0)

# Row 26 - سرمایه (Capital)
$ws.Range("D26").Value = 34888
$ws.Range("E26").Value = 27511
$ws.Range("F26").Value = 15609
$ws.Range("G26").Value = 13376
$ws.Range("H26").Value = 10001

# Row 27 - سود هر سهم بر اساس آخرین سرمایه : unchanged (0)
